$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("D2").Value = 0.16500606
$ws.Range("D3").Value = 0.34770248
